$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. MobileB sheet updates
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("MobileB")

# New header "Country" in K1
$ws.Range("K1").Value = "Country"

# Row 2 updates: phone number now a text value, card number split into
# last-4 (H2), exp month (I2) / exp year (I3), cvv stays 123, country added
$ws.Range("B2").Value = "'123456789"
$ws.Range("G2").Value = "'123456"
$ws.Range("H2").Value = "'5454"
$ws.Range("I2").Value = "'12"
$ws.Range("J2").Value = "'123"
$ws.Range("K2").Value = "India"

# New rows 3-5 with repeated masked card number fragments
$ws.Range("H3").Value = "'5454"
$ws.Range("I3").Value = "'28"
$ws.Range("H4").Value = "'5454"
$ws.Range("H5").Value = "'5454"

# Column width adjustments to match the new layout
$ws.Columns.Item(1).ColumnWidth = 8.85546875
$ws.Columns.Item(2).ColumnWidth = 9
$ws.Columns.Item(8).ColumnWidth = 8.43
$ws.Columns.Item(9).ColumnWidth = 17.28515625

# Update selection to J3
$ws.Range("J3").Select()

# Keep MobileB as the active/selected sheet
$ws.Activate()

# ------------------------------------------------------------------
# 2. New "Sheet2" with a second mobile-booking sample record
# ------------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $ws)
$new.Name = "Sheet2"

$new.Range("A1").Value = "FullName"
$new.Range("B1").Value = "PhoneN"
$new.Range("C1").Value = "Email"
$new.Range("D1").Value = "Address"
$new.Range("E1").Value = "Country"
$new.Range("F1").Value = "State"
$new.Range("G1").Value = "City"
$new.Range("H1").Value = "ZipCode"
$new.Range("I1").Value = "CardNum"
$new.Range("J1").Value = "ExpD"
$new.Range("K1").Value = "CCode"

$new.Range("A2").Value = "Dk"
$new.Range("B2").Value = "'7883664742"
$new.Range("C2").Value = "dinesh.kanna@igtsolutions.com"
$new.Range("D2").Value = "Chennai"
$new.Range("E2").Value = "India"
$new.Range("F2").Value = "Tamil Nadu"
$new.Range("G2").Value = "Chennai"
$new.Range("H2").Value = "'600064"
$new.Range("I2").Value = "'4012 8888 8888 1881"
$new.Range("J2").NumberFormat = "mmm-yy"
$new.Range("J2").Value = "11/22"
$new.Range("K2").Value = "'222"

$new.Columns.Item(2).ColumnWidth = 11
$new.Columns.Item(3).ColumnWidth = 30.140625
$new.Columns.Item(6).ColumnWidth = 11
$new.Columns.Item(9).ColumnWidth = 18.5703125

$new.Range("K3").Select()

# Re-activate MobileB so it remains the tab shown when the workbook opens
$ws.Activate()
